$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 6 (A6/B6 switch to bordered style s=6, C6/D6/E6 switch to bordered style s=7) ---
$ws.Range("A4").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("E4").Copy()
$ws.Range("E6").PasteSpecial(-4122)

# --- Row 7 (style set "4"/"5", no border) ---
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E7").PasteSpecial(-4122)

# --- Row 8 (style set "6"/"7", bordered) ---
$ws.Range("A4").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("E4").Copy()
$ws.Range("E8").PasteSpecial(-4122)

# --- Row 9 (style set "4"/"5", no border) ---
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E9").PasteSpecial(-4122)

# --- Row 10 (style set "6"/"7", bordered) ---
$ws.Range("A4").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("E4").Copy()
$ws.Range("E10").PasteSpecial(-4122)

# --- Row 11 (style set "4"/"5", no border) ---
$ws.Range("A2").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E11").PasteSpecial(-4122)

# --- Now fill in the cell values, in the exact order the strings were first
# authored, so the shared-string table indices line up with the target file. ---
$ws.Range("C7").Value = " There are many kinds of\ntreasures in the world."
$ws.Range("C8").Value = " We\'ve traveled far searching\nfor those treasures."
$ws.Range("A7").Value = "SCRIPT/T01P02A/um1201.ssb"
$ws.Range("D7").Value = " В мире есть много разных\nсокровищ."
$ws.Range("D8").Value = " Мы через многое прошли, пока\nискали эти сокровища."
$ws.Range("E7").Value = " Â íéñå åòóû íîïãï ñàèîúö\nòïëñïâéþ."
$ws.Range("E8").Value = " Íú œåñåè íîïãïå ðñïšìé, ðïëà\néòëàìé üóé òïëñïâéþà."
$ws.Range("A8").Value = "SCRIPT/T01P02A/um1203.ssb"
$ws.Range("C9").Value = " I had a really disgusting\ndrink there. It tasted terrible…"
$ws.Range("C10").Value = " I\'m not headed back to that\nshop again."
$ws.Range("A9").Value = "SCRIPT/P01P04A/um1404.ssb"
$ws.Range("D9").Value = " Мне тут подали такой\nотвратительный напиток. На вкус был\nужасен..."
$ws.Range("D10").Value = " Я больше не стану ничего\nзаказывать."
$ws.Range("E9").Value = " Íîå óôó ðïäàìé óàëïê\nïóâñàóéóåìûîúê îàðéóïë. Îà âëôò áúì\nôçàòåî..."
$ws.Range("E10").Value = " Ÿ áïìûšå îå òóàîô îéœåãï\nèàëàèúâàóû."
$ws.Range("A10").Value = " SCRIPT/P01P04A/um1504.ssb"
$ws.Range("C11").Value = " It would be good to see [CS:N]Grovyle[CR]\nfinally captured."
$ws.Range("A11").Value = "SCRIPT/G01P03A/um1613.ssb"
$ws.Range("D11").Value = " Будет здорово, если [CS:N]Гровайла[CR]\nнаконец поймают."
$ws.Range("E11").Value = " Áôäåó èäïñïâï, åòìé [CS:N]Ãñïâàêìà[CR]\nîàëïîåø ðïêíàýó."

# Numeric "row number in scripts" column
$ws.Range("B7").Value = 136
$ws.Range("B8").Value = 139
$ws.Range("B9").Value = 114
$ws.Range("B10").Value = 117
$ws.Range("B11").Value = 95

# --- Row heights matching the wrapped-text autosize in the source file ---
$ws.Rows.Item(7).RowHeight = 43.2
$ws.Rows.Item(8).RowHeight = 43.2
$ws.Rows.Item(9).RowHeight = 43.2
$ws.Rows.Item(10).RowHeight = 57.6
$ws.Rows.Item(11).RowHeight = 43.2

# --- Sheet view: scrolled down, new active cell selected ---
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("C15").Select()
